$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4166.594010438043
$ws.Range("C3").Value = 4166.594010438043
$ws.Range("C4").Value = 4163.662193025788
$ws.Range("C5").Value = 4163.662193025788
$ws.Range("C6").Value = 4143.202917589537
$ws.Range("C7").Value = 3962.055252700118
$ws.Range("C8").Value = 3962.055252700118
$ws.Range("C9").Value = 3957.918244635269
$ws.Range("C10").Value = 3898.995416859769
$ws.Range("C11").Value = 3898.995416859769
$ws.Range("C12").Value = 3898.995416859769
